# Applies the data-refresh update described in the commit diff
# (static market-price / profit values refreshed on several leve-profit sheets).
$wb = $excel.ActiveWorkbook

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 799.5599999999999
$ws.Range("I32").Value = 641.6951
$ws.Range("J32").Value = 1518.7222
$ws.Range("K32").Value = 641.6951
$ws.Range("L32").Value = 1518.7222
$ws.Range("M32").Value = -354.6951
$ws.Range("N32").Value = -2092.7222
$ws.Range("H122").Value = 1922.3864
$ws.Range("I122").Value = 1434.2059
$ws.Range("K122").Value = 4302.6177
$ws.Range("M122").Value = -1852.6177

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2509.8096
$ws.Range("I5").Value = 609.63635
$ws.Range("J5").Value = 4600
$ws.Range("K5").Value = 609.63635
$ws.Range("L5").Value = 4600
$ws.Range("M5").Value = -496.63635
$ws.Range("N5").Value = -4826

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1640.43
$ws.Range("I31").Value = 851.6286
$ws.Range("J31").Value = 2065.1692
$ws.Range("K31").Value = 851.6286
$ws.Range("L31").Value = 2065.1692
$ws.Range("M31").Value = -556.6286
$ws.Range("N31").Value = -2655.1692
$ws.Range("H34").Value = 1640.43
$ws.Range("I34").Value = 851.6286
$ws.Range("J34").Value = 2065.1692
$ws.Range("K34").Value = 851.6286
$ws.Range("L34").Value = 2065.1692
$ws.Range("M34").Value = -649.6286
$ws.Range("N34").Value = -2469.1692
$ws.Range("H86").Value = 2937.2307
$ws.Range("I86").Value = 3284.8823
$ws.Range("J86").Value = 2280.5557
$ws.Range("K86").Value = 3284.8823
$ws.Range("L86").Value = 2280.5557
$ws.Range("M86").Value = -2161.8823
$ws.Range("N86").Value = -4526.5557
$ws.Range("H89").Value = 2937.2307
$ws.Range("I89").Value = 3284.8823
$ws.Range("J89").Value = 2280.5557
$ws.Range("K89").Value = 16424.4115
$ws.Range("L89").Value = 11402.7785
$ws.Range("M89").Value = -10808.4115
$ws.Range("N89").Value = -22634.7785
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H132").Value = 41674664
$ws.Range("I132").Value = 76932630
$ws.Range("J132").Value = 6151.636
$ws.Range("K132").Value = 230797890
$ws.Range("L132").Value = 18454.908
$ws.Range("M132").Value = -230795360
$ws.Range("N132").Value = -23514.908

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1379.5333
$ws.Range("I122").Value = 1268.5
$ws.Range("J122").Value = 1601.6
$ws.Range("K122").Value = 3805.5
$ws.Range("L122").Value = 4804.799999999999
$ws.Range("M122").Value = -1355.5
$ws.Range("N122").Value = -9704.799999999999
$ws.Range("H123").Value = 10326
$ws.Range("J123").Value = 10326
$ws.Range("L123").Value = 10326
$ws.Range("N123").Value = -15226
$ws.Range("H126").Value = 1243.0385
$ws.Range("I126").Value = 1180.7333
$ws.Range("J126").Value = 1328
$ws.Range("K126").Value = 3542.199900000001
$ws.Range("L126").Value = 3984
$ws.Range("M126").Value = -1072.199900000001
$ws.Range("N126").Value = -8924
$ws.Range("H132").Value = 4218.4106
$ws.Range("I132").Value = 4947.0557
$ws.Range("J132").Value = 2906.85
$ws.Range("K132").Value = 14841.1671
$ws.Range("L132").Value = 8720.549999999999
$ws.Range("M132").Value = -12311.1671
$ws.Range("N132").Value = -13780.55

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 25560
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 25560
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 25560
$ws.Range("N119").Value = -35236
$ws.Range("H120").Value = 42710
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 42710
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 42710
$ws.Range("N120").Value = -52386
$ws.Range("H121").Value = 37420
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 37420
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 37420
$ws.Range("N121").Value = -40914
$ws.Range("H122").Value = 11168301
$ws.Range("I122").Value = 16751213
$ws.Range("J122").Value = 2476.8333
$ws.Range("K122").Value = 50253639
$ws.Range("L122").Value = 7430.499899999999
$ws.Range("M122").Value = -50251189
$ws.Range("N122").Value = -12330.4999
$ws.Range("H123").Value = 30214.5
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 30214.5
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 30214.5
$ws.Range("N123").Value = -40014.5
$ws.Range("H124").Value = 29429
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 29429
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 29429
$ws.Range("N124").Value = -39249
$ws.Range("H125").Value = 47857.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 47857.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 47857.5
$ws.Range("N125").Value = -57697.5
$ws.Range("H126").Value = 1046.1333
$ws.Range("I126").Value = 1007.6667
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 3023.0001
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -553.0001000000002
$ws.Range("N126").Value = -8540
$ws.Range("H127").Value = 28714.5
$ws.Range("I127").Value = 20000
$ws.Range("J127").Value = 37429
$ws.Range("K127").Value = 20000
$ws.Range("L127").Value = 37429
$ws.Range("M127").Value = -15040
$ws.Range("N127").Value = -47349
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("H129").Value = 36357.25
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 36357.25
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 36357.25
$ws.Range("N129").Value = -46357.25
$ws.Range("H130").Value = 36398
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 36398
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 36398
$ws.Range("N130").Value = -46438
$ws.Range("H131").Value = 42857.5
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 42857.5
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 42857.5
$ws.Range("N131").Value = -52937.5
$ws.Range("H132").Value = 7465683
$ws.Range("I132").Value = 9807033
$ws.Range("J132").Value = 2630.125
$ws.Range("K132").Value = 29421099
$ws.Range("L132").Value = 7890.375
$ws.Range("M132").Value = -29418569
$ws.Range("N132").Value = -12950.375
$ws.Range("H133").Value = 44536.25
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 44536.25
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 44536.25
$ws.Range("N133").Value = -54656.25
$ws.Range("H135").Value = 49715
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 49715
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 49715
$ws.Range("N135").Value = -59855
$ws.Range("H136").Value = 3089.2964
$ws.Range("I136").Value = 4998.1304
$ws.Range("J136").Value = 1673.0646
$ws.Range("K136").Value = 14994.3912
$ws.Range("L136").Value = 5019.1938
$ws.Range("M136").Value = -12444.3912
$ws.Range("N136").Value = -10119.1938
$ws.Range("H137").Value = 44715
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 44715
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 44715
$ws.Range("N137").Value = -54915
$ws.Range("H138").Value = 44429
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 44429
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 44429
$ws.Range("N138").Value = -54709
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 38685.8
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 38685.8
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 38685.8
$ws.Range("N140").Value = -49045.8
$ws.Range("H141").Value = 38796.43
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 38796.43
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 38796.43
$ws.Range("N141").Value = -49156.43

Write-Output "Applied Ramuh_Profits data refresh across ARM, BSM, CRP, GSM, WVR sheets."
